$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.274.67"
$ws.Range("E2").Value = "  +4.91%  "

# Row 3
$ws.Range("D3").Value = "2.615.66"
$ws.Range("E3").Value = "  +4.98%  "

# Row 5
$ws.Range("D5").Value = "'606.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.98%  "

# Row 6
$ws.Range("D6").Value = "'180.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.56%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +1.80%  "

# Row 9
$ws.Range("D9").Value = "2.614.96"
$ws.Range("E9").Value = "  +4.98%  "

# Row 10
$ws.Range("D10").Value = "'0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.82%  "

# Row 11
$ws.Range("D11").Value = "'0.165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12
$ws.Range("D12").Value = "'0.350"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.94%  "

# Row 13
$ws.Range("E13").Value = "  +1.66%  "

# Row 14
$ws.Range("D14").Value = "3.070.49"
$ws.Range("E14").Value = "  +4.31%  "

# Row 15
$ws.Range("D15").Value = "'26.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.52%  "

# Row 16
$ws.Range("E16").Value = "  +7.79%  "

# Row 17
$ws.Range("D17").Value = "71.255.50"
$ws.Range("E17").Value = "  +4.42%  "

# Row 18
$ws.Range("D18").Value = "2.610.32"
$ws.Range("E18").Value = "  +4.30%  "

# Row 19
$ws.Range("D19").Value = "'7.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.75%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'378.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.09%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'11.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.88%  "

# Row 22
$ws.Range("E22").Value = "  +1.58%  "

# Row 23
$ws.Range("D23").Value = "'72.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.42%  "

# Row 24
$ws.Range("D24").Value = "'4.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.68%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("D26").Value = "'1.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.81%  "

# Row 27
$ws.Range("D27").Value = "'9.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.91%  "

# Row 28
$ws.Range("D28").Value = "2.750.41"
$ws.Range("E28").Value = "  +6.51%  "

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0956"
$ws.Range("E30").Value = "  +6.05%  "

# Row 31
$ws.Range("D31").Value = "'533.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.33%  "

# Row 32
$ws.Range("E32").Value = "  +4.21%  "

# Row 33
$ws.Range("D33").Value = "'1.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.90%  "

# Row 34
$ws.Range("E34").Value = "  +4.16%  "

# Row 35
$ws.Range("E35").Value = "  +0.10%  "

# Row 36
$ws.Range("D36").Value = "'164.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.81%  "

# Row 37
$ws.Range("D37").Value = "'0.120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.55%  "

# Row 38
$ws.Range("D38").Value = "'19.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.82%  "

# Row 39
$ws.Range("E39").Value = "  +1.74%  "

# Row 40
$ws.Range("D40").Value = "'1.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.95%  "

# Row 41
$ws.Range("D41").Value = "'1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.94%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.28%  "

# Row 43
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("E44").Value = "  +5.38%  "

# Row 45
$ws.Range("D45").Value = "'0.331"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.75%  "

# Row 46
$ws.Range("D46").Value = "'40.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.58%  "

# Row 47
$ws.Range("D47").Value = "'154.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.24%  "

# Row 48
$ws.Range("E48").Value = "  +3.79%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.534"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.96%  "

# Row 50
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "'1.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.43%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0755"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.22%  "

